$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("Data")
$codeWs = $wb.Worksheets.Item("Codebook")

# Update the "Data" sheet: rename D/E headers and replace eye-color/waist
# columns with inseam (numeric) / hair color (categorical) data.
$dataWs.Range("D1").Value = "Inseam"
$dataWs.Range("E1").Value = "Hair Color"

$dataValues = @(
    @(81, "black"),
    @(74, "blond"),
    @(55, "brown"),
    @(91, "l brown"),
    @(95, "other"),
    @(89, "black"),
    @(68, "other"),
    @(62, "blond"),
    @(73, "brown"),
    @(49, "other"),
    @(53, "other"),
    @(50, "l brown"),
    @(52, "d brown"),
    @(65, "white")
)

for ($i = 0; $i -lt $dataValues.Count; $i++) {
    $row = $i + 2
    $dataWs.Cells.Item($row, 4).Value = $dataValues[$i][0]
    $dataWs.Cells.Item($row, 5).Value = $dataValues[$i][1]
}

# Update the "Codebook" sheet: remove the now-obsolete "Eye Color" and
# "Waist" rows (rows 5 and 6).
$codeWs.Rows.Item(6).Delete()
$codeWs.Rows.Item(5).Delete()

$dataWs.Range("D22").Select()
$codeWs.Range("A4").Select()
$dataWs.Activate()
